$d = $word.ActiveDocument

# Locate the paragraph that holds "LOM3257: Mecânica Clássica (Requisito fraco)" —
# the trailing footer block (blank line, "Ver no Jupiter..." line, and the
# "© 2020 ..." copyright line) that follows it is what the Jekyll site rebuild
# dropped, so find its bounds dynamically rather than hard-coding indices.
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOM3257*") {
        $startIndex = $i + 1
    }
    if ($t -like "*Original theme under Creative Commons Attribution*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $victim = $d.Range($rangeStart, $rangeEnd)
    $victim.Delete()
}
